$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The edit inserts one new weekly observation at the top of the data block
# (row 74), which pushes every existing observation down by one row; the
# observation that used to be the very last one (row 185) becomes the new
# last row (186). Columns A,B,C,E,F,G,H,I,J,K,T are constant across the whole
# block, so only D (Fecha) and L..S (Calidad..Precio $/Kg) actually move.
# ---------------------------------------------------------------------------

$firstRow = 74
$lastRow  = 185
$newLastRow = 186

# --- 1. cache the original values for the columns that move (D, L..S) -----
$D = @{}
$L = @{}
$M = @{}
$N = @{}
$O = @{}
$P = @{}
$Q = @{}
$R = @{}
$S = @{}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $D[$r] = $ws.Cells.Item($r, 4).Value2
    $L[$r] = $ws.Cells.Item($r, 12).Value2
    $M[$r] = $ws.Cells.Item($r, 13).Value2
    $N[$r] = $ws.Cells.Item($r, 14).Value2
    $O[$r] = $ws.Cells.Item($r, 15).Value2
    $P[$r] = $ws.Cells.Item($r, 16).Value2
    $Q[$r] = $ws.Cells.Item($r, 17).Value2
    $R[$r] = $ws.Cells.Item($r, 18).Value2
    $S[$r] = $ws.Cells.Item($r, 19).Value2
}

# --- 2. build the brand new row 186 out of the (still intact) static ------
#        columns of row 185, plus the D/L..S values that used to live in
#        row 185 (the old "last" observation moves to the new last row).
for ($col = 1; $col -le 20; $col++) {
    $ws.Cells.Item($newLastRow, $col).Value = $ws.Cells.Item($lastRow, $col).Value2
}
$ws.Range("D$newLastRow").NumberFormat = $ws.Range("D$lastRow").NumberFormat

# --- 3. shift every row down by one: new row i gets what used to be in ----
#        row i-1 (for i = lastRow down to firstRow+1), using the cached
#        original values so we never read data that has already been
#        overwritten.
for ($r = $lastRow; $r -ge ($firstRow + 1); $r--) {
    $src = $r - 1
    $ws.Cells.Item($r, 4).Value  = $D[$src]
    $ws.Cells.Item($r, 12).Value = $L[$src]
    $ws.Cells.Item($r, 13).Value = $M[$src]
    $ws.Cells.Item($r, 14).Value = $N[$src]
    $ws.Cells.Item($r, 15).Value = $O[$src]
    $ws.Cells.Item($r, 16).Value = $P[$src]
    $ws.Cells.Item($r, 17).Value = $Q[$src]
    $ws.Cells.Item($r, 18).Value = $R[$src]
    $ws.Cells.Item($r, 19).Value = $S[$src]
}

# --- 4. row 74 becomes the brand new observation ---------------------------
$ws.Cells.Item($firstRow, 4).Value  = 44579
$ws.Cells.Item($firstRow, 12).Value = "Primera"
$ws.Cells.Item($firstRow, 13).Value = 1000
$ws.Cells.Item($firstRow, 14).Value = 8500
$ws.Cells.Item($firstRow, 15).Value = 9000
$ws.Cells.Item($firstRow, 16).Value = 8750
$ws.Cells.Item($firstRow, 17).Value = "`$/caja 7 kilos"
$ws.Cells.Item($firstRow, 18).Value = "Región de La Araucanía"
$ws.Cells.Item($firstRow, 19).Value = 1250

Write-Host "edit complete"
